$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values that look like plain numbers (e.g. "398.44") must be written into
# cells pre-formatted as Text, otherwise Excel auto-converts the typed text into
# a real number, which would change both the stored type and the on-screen text
# (e.g. trailing zeros in "1.00"/"22.30" would be lost). These "Price" column
# cells are plain inline text in the source workbook, so preserve that.
$textCells = @(
    "D5", "D6", "D10", "D11", "D14", "D18", "D21", "D22", "D23", "D24"
    "D25", "D26", "D27", "D28", "D29", "D34", "D36", "D39", "D40", "D41"
    "D44", "D45", "D46", "D47", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '56.358.74'
$ws.Range('E2').Value = '  +9.90%  '
$ws.Range('D3').Value = '3.224.63'
$ws.Range('E3').Value = '  +4.20%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '398.44'
$ws.Range('E5').Value = '  +2.58%  '
$ws.Range('D6').Value = '111.04'
$ws.Range('E6').Value = '  +7.12%  '
$ws.Range('E9').Value = '  +6.15%  '
$ws.Range('D10').Value = '39.28'
$ws.Range('E10').Value = '  +6.53%  '
$ws.Range('D11').Value = '0.0922'
$ws.Range('E11').Value = '  +7.70%  '
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('D13').Value = '3.735.02'
$ws.Range('E13').Value = '  +4.25%  '
$ws.Range('D14').Value = '8.08'
$ws.Range('E14').Value = '  +4.52%  '
$ws.Range('E15').Value = '  +3.05%  '
$ws.Range('D16').Value = '3.215.60'
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('E17').Value = '  +4.97%  '
$ws.Range('D18').Value = '10.89'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '56.240.69'
$ws.Range('E19').Value = '  +9.45%  '
$ws.Range('E20').Value = '  +3.33%  '
$ws.Range('D21').Value = '0.0000104'
$ws.Range('E21').Value = '  +7.47%  '
$ws.Range('D22').Value = '13.02'
$ws.Range('E22').Value = '  +4.62%  '
$ws.Range('D23').Value = '296.86'
$ws.Range('E23').Value = '  +11.81%  '
$ws.Range('D24').Value = '75.76'
$ws.Range('E24').Value = '  +8.24%  '
$ws.Range('D25').Value = '3.23'
$ws.Range('E25').Value = '  +1.86%  '
$ws.Range('D26').Value = '8.16'
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').Value = '28.05'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('D28').Value = '7.42'
$ws.Range('E28').Value = '  +2.66%  '
$ws.Range('D29').Value = '0.172'
$ws.Range('E29').Value = '  +4.21%  '
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('E31').Value = '  +3.89%  '
$ws.Range('E32').Value = '  +6.52%  '
$ws.Range('E33').Value = '  +4.03%  '
$ws.Range('D34').Value = '36.61'
$ws.Range('E35').Value = '  +3.25%  '
$ws.Range('D36').Value = '51.32'
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('E37').Value = '  +25.13%  '
$ws.Range('E38').Value = '  +4.05%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '134.78'
$ws.Range('E40').Value = '  +3.29%  '
$ws.Range('D41').Value = '17.37'
$ws.Range('E41').Value = '  +4.43%  '
$ws.Range('E42').Value = '  +3.62%  '
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('D44').Value = '0.120'
$ws.Range('E44').Value = '  +3.33%  '
$ws.Range('D45').Value = '0.283'
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '22.30'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').Value = '2.20'
$ws.Range('E47').Value = '  +54.28%  '
$ws.Range('D48').Value = '2.129.56'
$ws.Range('E48').Value = '  +2.80%  '
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').Value = '2.44'
$ws.Range('E50').Value = '  -2.98%  '
$ws.Range('E51').Value = '  +11.57%  '
